$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab
$ws.Name = "Rajat Patidar"

# Insert a new column A ("matchNo") - shifts existing columns B..M right
$ws.Columns.Item(1).Insert()

# Insert 3 new blank rows above the existing data row (old row 2),
# pushing it down to row 5, making room for 3 new match rows (2,3,4)
$ws.Rows.Item(2).Resize(3).Insert()

# Force the numeric-looking new cells to stay text (matches source data
# which is stored as text, like the rest of the sheet)
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("E2:I4").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2 - match 22nd
$ws.Range("A2").Value = "22nd"
$ws.Range("B2").Value = "Royal Challengers Bangalore"
$ws.Range("C2").Value = "Rajat Patidar"
$ws.Range("D2").Value = "c Smith b Patel"
$ws.Range("E2").Value = "31"
$ws.Range("F2").Value = "22"
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "2"
$ws.Range("I2").Value = "140.90"
$ws.Range("J2").Value = "Delhi Capitals"
$ws.Range("K2").Value = "Ahmedabad"
$ws.Range("L2").Value = "April 27"
$ws.Range("M2").Value = "RCB won by 1 run"

# Row 3 - match 10th
$ws.Range("A3").Value = "10th"
$ws.Range("B3").Value = "Royal Challengers Bangalore"
$ws.Range("C3").Value = "Rajat Patidar"
$ws.Range("D3").Value = "b Varun"
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "2"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "50.00"
$ws.Range("J3").Value = "Kolkata Knight Riders"
$ws.Range("K3").Value = "Chennai"
$ws.Range("L3").Value = "April 18"
$ws.Range("M3").Value = "RCB won by 38 runs"

# Row 4 - match 26th
$ws.Range("A4").Value = "26th"
$ws.Range("B4").Value = "Royal Challengers Bangalore"
$ws.Range("C4").Value = "Rajat Patidar"
$ws.Range("D4").Value = "c Pooran b Jordan"
$ws.Range("E4").Value = "31"
$ws.Range("F4").Value = "30"
$ws.Range("G4").Value = "2"
$ws.Range("H4").Value = "1"
$ws.Range("I4").Value = "103.33"
$ws.Range("J4").Value = "Punjab Kings"
$ws.Range("K4").Value = "Ahmedabad"
$ws.Range("L4").Value = "April 30"
$ws.Range("M4").Value = "Punjab Kings won by 34 runs"

# Row 5 already holds the original match (1st) data in columns B..M,
# shifted down/right by the column/row inserts above and still correctly
# typed as text there. Only the brand-new column A cell needs a value.
$ws.Range("A5").Value = "1st"
